# TS 4.6 Jatai Sanskrit Corrections – update the "Observed till" date from
# "30th April 2024" to "31st May 2024".
#
# The original text lives in a single run. The author's actual edit (per
# the target OOXML diff) typed the replacement over the old text in stages,
# so Word ended up keeping the identical run formatting (bold, bold-complex-
# script, size 32/32cs, single underline) but split across five runs:
#   "3" | "1st" | " " | "May" | " 2024"
#
# Plain Find/Replace (or Range.Text assignment) would just merge everything
# back into a single run because the formatting is identical, so instead we
# locate the exact Range of the old text and replace its contents with a
# small WordprocessingML fragment (via Range.InsertXML) that reproduces the
# five-run split exactly.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("30th April 2024", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the text '30th April 2024' to replace."
}

$rPr = '<w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr>'

$runsXml = (
    "<w:r>$rPr<w:t>3</w:t></w:r>" +
    "<w:r>$rPr<w:t>1st</w:t></w:r>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:r>$rPr<w:t>May</w:t></w:r>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> 2024</w:t></w:r>"
)

$packageXml = (
    '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" ' +
    'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    "<w:body><w:p>$runsXml</w:p></w:body>" +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
)

$rng.InsertXML($packageXml)
